$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric (e.g. "168.10", "36.482.08").
# Force the whole price column to Text format BEFORE writing so Excel does not
# silently coerce these into numbers (which would drop significant trailing
# zeros / thousand-separator dots), matching the original inlineStr storage.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '36.483.37'
$ws.Range("E2").Value = '  +3.09%  '
$ws.Range("D3").Value = '1.924.08'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '249.71'
$ws.Range("E5").Value = '  +1.29%  '
$ws.Range("E6").Value = '  +0.45%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '43.94'
$ws.Range("E8").Value = '  +1.50%  '
$ws.Range("D9").Value = '57.87'
$ws.Range("E9").Value = '  +7.57%  '
$ws.Range("E10").Value = '  +3.28%  '
$ws.Range("E11").Value = '  +3.11%  '
$ws.Range("D12").Value = '0.0996'
$ws.Range("E12").Value = '  +2.66%  '
$ws.Range("D13").Value = '14.61'
$ws.Range("E13").Value = '  +9.73%  '
$ws.Range("D14").Value = '0.798'
$ws.Range("E14").Value = '  +5.12%  '
$ws.Range("D15").Value = '2.201.80'
$ws.Range("E15").Value = '  +1.67%  '
$ws.Range("D16").Value = '5.13'
$ws.Range("E16").Value = '  +4.81%  '
$ws.Range("D17").Value = '1.928.06'
$ws.Range("E17").Value = '  +2.12%  '
$ws.Range("D18").Value = '36.459.73'
$ws.Range("E18").Value = '  +2.81%  '
$ws.Range("D19").Value = '74.59'
$ws.Range("E19").Value = '  +1.89%  '
$ws.Range("D20").Value = '0.0₃0845'
$ws.Range("E20").Value = '  +2.47%  '
$ws.Range("D21").Value = '253.89'
$ws.Range("E21").Value = '  +3.60%  '
$ws.Range("D22").Value = '13.27'
$ws.Range("E22").Value = '  +3.51%  '
$ws.Range("D23").Value = '5.23'
$ws.Range("E23").Value = '  +5.52%  '
$ws.Range("E24").Value = '  +1.16%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").Value = '2.27'
$ws.Range("E26").Value = '  +5.33%  '
$ws.Range("D27").Value = '168.10'
$ws.Range("E27").Value = '  +1.16%  '
$ws.Range("D28").Value = '8.83'
$ws.Range("E28").Value = '  +3.68%  '
$ws.Range("D29").Value = '18.91'
$ws.Range("E29").Value = '  +2.93%  '
$ws.Range("E30").Value = '  +1.49%  '
$ws.Range("D31").Value = '4.55'
$ws.Range("E31").Value = '  +6.98%  '
$ws.Range("E32").Value = '  +4.45%  '
$ws.Range("B33").Value = 'WEMIXToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D33").Value = '1.95'
$ws.Range("E33").Value = '  +5.44%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '4.34'
$ws.Range("E34").Value = '  +3.97%  '
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  +21.73%  '
$ws.Range("D37").Value = '1.49'
$ws.Range("E37").Value = '  -14.43%  '
$ws.Range("E38").Value = '  +2.20%  '
$ws.Range("D39").Value = '2.02'
$ws.Range("E39").Value = '  +2.60%  '
$ws.Range("D40").Value = '104.96'
$ws.Range("E40").Value = '  +8.27%  '
$ws.Range("D41").Value = '0.0229'
$ws.Range("E41").Value = '  +3.99%  '
$ws.Range("D42").Value = '17.21'
$ws.Range("E42").Value = '  -0.61%  '
$ws.Range("D43").Value = '15.10'
$ws.Range("E43").Value = '  +22.40%  '
$ws.Range("E44").Value = '  +3.25%  '
$ws.Range("D45").Value = '1.345.78'
$ws.Range("E45").Value = '  +3.54%  '
$ws.Range("E46").Value = '  +3.00%  '
$ws.Range("E47").Value = '  +1.86%  '
$ws.Range("E48").Value = '  +1.84%  '
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("E50").Value = '  +2.88%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.101.98'
$ws.Range("E51").Value = '  +1.37%  '
